# Minor Bugfixing Planung A, B, D 18.11.2021
# Updates the "Sprint-Backlog" sheet of the Team A SCRUM workbook:
#  - mark rows 7 & 8 as "Done"
#  - rename / re-estimate the "Time-Limit" task (row 9)
#  - clear out row 10 (task moved / removed)
#  - row 11 moves from sprint-count 3 to 2 and is no longer "Überplant"
#  - row 12 gets a new backlog item describing the new game's UI
#  - selection / scroll position is left where the user last worked (B12)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint-Backlog")
$ws.Activate()

# Row 7: implementation of the human TicTacToe player is finished
$ws.Range("E7").Value = "Done"

# Row 8: implementation of the computer TicTacToe player is finished
$ws.Range("E8").Value = "Done"

# Row 9: the "Time-Limit" item became "Time-Limit fertig machen" and its
# estimate changed from 180min to 120min
$ws.Range("B9").Value = "Time-Limit fertig machen"
$ws.Range("D9").Value = "120min"

# Row 10: this item was removed from the sprint backlog entirely
$ws.Range("A10:E10").ClearContents()

# Row 11: "Spezifikation neues Spiel" moves back from sprint 3 to sprint 2
# and is no longer flagged as "Überplant"
$ws.Range("A11").Value = 2
$ws.Range("E11").ClearContents()

# Row 12: new backlog item for the new game's UI / interface concept
$ws.Range("A12").Value = 2
$ws.Range("C12").Value = "Darstellung für neu spezifiziertes Spiel. Painter-Objekt."
$ws.Range("B12").Value = "Darstellung neues Spiel inkl. Schnittstellenkonzept"
$ws.Range("D12").Value = "180min"

# Leave the view/selection where the user ended up editing
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B12").Select()
